$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for each data row (2-135).
# All of them are being updated from 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C135").Value = 45175
